# Challenge_1_50days_tracker.xlsx
# Commit: "226. Invert Binary Tree(recursion)"
#
# Adds a new "Binary Tree (basic)" section header (row 55) and a new
# data row (row 57) for LeetCode #226 "Invert Binary Tree", leaving a
# blank separator row (56) just like the rest of the sheet's sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New section header row (row 55) -------------------------------------
# Clone the formatting of an existing "module header" row (row 51, style
# index 7 - same blue-ish header style used from row 42 onward) onto the
# new header cell, then overwrite its text.
$ws.Range("B51").Copy($ws.Range("B55"))
$ws.Range("B55").Value = "Binary Tree (basic)"
$ws.Rows.Item(55).RowHeight = 35

# --- New question row (row 57) --------------------------------------------
$ws.Range("A57").Value = 226
$ws.Range("B57").Value = "Invert Binary Tree"
$ws.Range("C57").Value = "Easy"
$ws.Range("D57").Value = "Binary Tree,recurson"
$ws.Range("E57").Value = 45733
$ws.Range("E57").NumberFormat = "dd\-mmm\-yy"

# --- Keep the view pointed at the new bottom of the sheet ------------------
$ws.Range("E58").Select() | Out-Null
